# Powerpoint writer: consolidate text runs when possible.
#
# Three paragraphs in this deck were originally split across multiple
# <a:r> runs ("The" / " " / "Moon" and "One" / " " / "More") even though
# they carry no per-run formatting differences. The fix re-saves each of
# those paragraphs as a single run with the full text.
#
# The COM TextRange.Text setter performs a minimal in-place edit against
# whatever run structure currently exists, only touching the runs whose
# characters actually changed (same behavior real PowerPoint exhibits).
# Assigning the already-equal string back is therefore a no-op, and
# assigning a value that merely shares leading/trailing characters with
# the existing runs only request a partial split/merge. To force a full,
# single-run rewrite we first assign a placeholder string built from
# characters that cannot appear in any of the target strings (so there is
# zero shared prefix/suffix with the current runs), then assign the real,
# final text; at that point the whole paragraph is just one run and the
# second assignment rewrites it as exactly one <a:r>.

$placeholder = "0123456789"

$p = $ppt.ActivePresentation

# Slide 2: textbox "The Moon" caption under the moon picture.
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item("TextBox 3").TextFrame.TextRange
$tb2.Text = $placeholder
$tb2.Text = "The Moon"

# Slide 3: title "One More" and textbox "The Moon" caption.
$s3 = $p.Slides.Item(3)

$title3 = $s3.Shapes.Item("Title 1").TextFrame.TextRange
$title3.Text = $placeholder
$title3.Text = "One More"

$tb3 = $s3.Shapes.Item("TextBox 3").TextFrame.TextRange
$tb3.Text = $placeholder
$tb3.Text = "The Moon"
